$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rngB = $ws.Range("B2:B51")
$rngC = $ws.Range("C2:C51")
$rngD = $ws.Range("D2:D51")
$rngE = $ws.Range("E2:E51")

$arrB = New-Object "object[,]" 50,1
$arrC = New-Object "object[,]" 50,1
$arrD = New-Object "object[,]" 50,1
$arrE = New-Object "object[,]" 50,1

$arrB[0,0] = "Bitcoin"
$arrC[0,0] = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$arrD[0,0] = "26.722.78"
$arrE[0,0] = "  +4.16%  "
$arrB[1,0] = "Ethereum"
$arrC[1,0] = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$arrD[1,0] = "1.873.99"
$arrE[1,0] = "  +3.58%  "
$arrB[2,0] = "TetherUSD"
$arrC[2,0] = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$arrD[2,0] = "0.9998"
$arrE[2,0] = "  -0.19%  "
$arrB[3,0] = "BNB"
$arrC[3,0] = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$arrD[3,0] = "282.33"
$arrE[3,0] = "  +2.07%  "
$arrB[4,0] = "USDC"
$arrC[4,0] = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$arrD[4,0] = "0.9996"
$arrE[4,0] = "  -0.20%  "
$arrB[5,0] = "XRP"
$arrC[5,0] = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$arrD[5,0] = "0.5159"
$arrE[5,0] = "  +2.92%  "
$arrB[6,0] = "Cardano"
$arrC[6,0] = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$arrD[6,0] = "0.3540"
$arrE[6,0] = "  +1.17%  "
$arrB[7,0] = "OKB"
$arrC[7,0] = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$arrD[7,0] = "45.31"
$arrE[7,0] = "  +3.40%  "
$arrB[8,0] = "Dogecoin"
$arrC[8,0] = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$arrD[8,0] = "0.07119"
$arrE[8,0] = "  +7.52%  "
$arrB[9,0] = "Solana"
$arrC[9,0] = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$arrD[9,0] = "20.19"
$arrE[9,0] = "  +0.66%  "
$arrB[10,0] = "Polygon"
$arrC[10,0] = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$arrD[10,0] = "0.8203"
$arrE[10,0] = "  -2.17%  "
$arrB[11,0] = "TRON"
$arrC[11,0] = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$arrD[11,0] = "0.07759"
$arrE[11,0] = "  -0.80%  "
$arrB[12,0] = "WrappedEther"
$arrC[12,0] = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$arrD[12,0] = "1.862.49"
$arrE[12,0] = "  +2.88%  "
$arrB[13,0] = "Polkadot"
$arrC[13,0] = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$arrD[13,0] = "5.162"
$arrE[13,0] = "  +2.55%  "
$arrB[14,0] = "Litecoin"
$arrC[14,0] = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$arrD[14,0] = "89.62"
$arrE[14,0] = "  +2.62%  "
$arrB[15,0] = "BinanceUSD"
$arrC[15,0] = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$arrD[15,0] = "0.9993"
$arrE[15,0] = "  -0.12%  "
$arrB[16,0] = "Avalanche"
$arrC[16,0] = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$arrD[16,0] = "14.43"
$arrE[16,0] = "  +4.09%  "
$arrB[17,0] = "ShibaInu"
$arrC[17,0] = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$arrD[17,0] = "0.000008177"
$arrE[17,0] = "  +3.01%  "
$arrB[18,0] = "Dai"
$arrC[18,0] = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$arrD[18,0] = "0.9989"
$arrE[18,0] = "  -0.25%  "
$arrB[19,0] = "WrappedBTC"
$arrC[19,0] = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$arrD[19,0] = "26.768.19"
$arrE[19,0] = "  +4.05%  "
$arrB[20,0] = "Uniswap"
$arrC[20,0] = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$arrD[20,0] = "4.782"
$arrE[20,0] = "  +1.57%  "
$arrB[21,0] = "Cosmos"
$arrC[21,0] = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$arrD[21,0] = "10.17"
$arrE[21,0] = "  +1.87%  "
$arrB[22,0] = "Chainlink"
$arrC[22,0] = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$arrD[22,0] = "6.254"
$arrE[22,0] = "  +3.27%  "
$arrB[23,0] = "LidoDAOToken"
$arrC[23,0] = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$arrD[23,0] = "2.422"
$arrE[23,0] = "  +16.09%  "
$arrB[24,0] = "Monero"
$arrC[24,0] = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$arrD[24,0] = "145.71"
$arrE[24,0] = "  +3.18%  "
$arrB[25,0] = "EthereumClassic"
$arrC[25,0] = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$arrD[25,0] = "17.44"
$arrE[25,0] = "  +3.53%  "
$arrB[26,0] = "Toncoin"
$arrC[26,0] = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$arrD[26,0] = "1.661"
$arrE[26,0] = "  +0.08%  "
$arrB[27,0] = "BitcoinCash"
$arrC[27,0] = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$arrD[27,0] = "111.01"
$arrE[27,0] = "  +2.30%  "
$arrB[28,0] = "InternetComputer(DFINITY)"
$arrC[28,0] = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$arrD[28,0] = "4.422"
$arrE[28,0] = "  +2.97%  "
$arrB[29,0] = "Filecoin"
$arrC[29,0] = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$arrD[29,0] = "4.363"
$arrE[29,0] = "  +3.99%  "
$arrB[30,0] = "Stellar"
$arrC[30,0] = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$arrD[30,0] = "0.08834"
$arrE[30,0] = "  +0.52%  "
$arrB[31,0] = "Hedera"
$arrC[31,0] = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$arrD[31,0] = "0.04904"
$arrE[31,0] = "  +1.62%  "
$arrB[32,0] = "ARBITRUM"
$arrC[32,0] = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$arrD[32,0] = "1.180"
$arrE[32,0] = "  +5.08%  "
$arrB[33,0] = "ImmutableX"
$arrC[33,0] = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$arrD[33,0] = "0.7491"
$arrE[33,0] = "  +1.53%  "
$arrB[34,0] = "HuobiToken"
$arrC[34,0] = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$arrD[34,0] = "2.862"
$arrE[34,0] = "  -0.29%  "
$arrB[35,0] = "MXToken"
$arrC[35,0] = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$arrD[35,0] = "3.279"
$arrE[35,0] = "  +7.93%  "
$arrB[36,0] = "RenderToken"
$arrC[36,0] = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$arrD[36,0] = "2.445"
$arrE[36,0] = "  +2.30%  "
$arrB[37,0] = "TheSandbox"
$arrC[37,0] = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$arrD[37,0] = "0.5310"
$arrE[37,0] = "  +2.37%  "
$arrB[38,0] = "VeChain"
$arrC[38,0] = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$arrD[38,0] = "0.01886"
$arrE[38,0] = "  +1.16%  "
$arrB[39,0] = "TrustWalletToken"
$arrC[39,0] = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$arrD[39,0] = "0.9770"
$arrE[39,0] = "  +0.54%  "
$arrB[40,0] = "Quant"
$arrC[40,0] = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$arrD[40,0] = "116.43"
$arrE[40,0] = "  +4.65%  "
$arrB[41,0] = "FraxShare"
$arrC[41,0] = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$arrD[41,0] = "6.307"
$arrE[41,0] = "  +1.70%  "
$arrB[42,0] = "Aptos"
$arrC[42,0] = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$arrD[42,0] = "8.175"
$arrE[42,0] = "  +0.68%  "
$arrB[43,0] = "PaxDollar"
$arrC[43,0] = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$arrD[43,0] = "0.9990"
$arrE[43,0] = "  -0.20%  "
$arrB[44,0] = "Decentraland"
$arrC[44,0] = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$arrD[44,0] = "0.4621"
$arrE[44,0] = "  +0.41%  "
$arrB[45,0] = "Algorand"
$arrC[45,0] = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$arrD[45,0] = "0.1365"
$arrE[45,0] = "  -1.01%  "
$arrB[46,0] = "EnergySwap"
$arrC[46,0] = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$arrD[46,0] = "9.518"
$arrE[46,0] = "  +3.28%  "
$arrB[47,0] = "Elrond"
$arrC[47,0] = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$arrD[47,0] = "36.67"
$arrE[47,0] = "  +2.90%  "
$arrB[48,0] = "NEARProtocol"
$arrC[48,0] = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$arrD[48,0] = "1.519"
$arrE[48,0] = "  +2.18%  "
$arrB[49,0] = "Cronos"
$arrC[49,0] = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$arrD[49,0] = "0.05930"
$arrE[49,0] = "  +1.34%  "

$rngB.NumberFormat = "@"
$rngC.NumberFormat = "@"
$rngD.NumberFormat = "@"
$rngE.NumberFormat = "@"

$rngB.Value = $arrB
$rngC.Value = $arrC
$rngD.Value = $arrD
$rngE.Value = $arrE

$rngB.NumberFormat = "General"
$rngC.NumberFormat = "General"
$rngD.NumberFormat = "General"
$rngE.NumberFormat = "General"

$rngB.Style = "Normal"
$rngC.Style = "Normal"
$rngD.Style = "Normal"
$rngE.Style = "Normal"
